$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the bold/centered/bordered formatting already used for the label
# column (A2:A16) on the three brand-new rows appended at the bottom.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# The "Gaussian-Quadrature" scheme (previously the last row, 16) now sits
# right after the "Ring Perpendicular" rows; three freshly-run spiral
# schemes are inserted after it, and the remaining rotation/hex-grid runs
# shift down underneath them, extending the table to row 19.

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$row10 = New-Object 'object[,]' 1,14
$row10[0,0] = 1.034373068099897
$row10[0,1] = 0.9546904888427717
$row10[0,2] = 0.9948365564448869
$row10[0,3] = 0.9842844965884121
$row10[0,4] = 1.034373068099897
$row10[0,5] = 0.9546904888427717
$row10[0,6] = 1.009099559340002
$row10[0,7] = 0.9803122053658803
$row10[0,8] = 1.008112316780786
$row10[0,9] = 0.9676901038902267
$row10[0,10] = 1.034373068099897
$row10[0,11] = 0.9747635226438294
$row10[0,12] = 0.9920461524939921
$row10[0,13] = 0.9916748494191079
$ws.Range("C10:P10").Value = $row10

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$row11 = New-Object 'object[,]' 1,14
$row11[0,0] = 1.037776720066717
$row11[0,1] = 0.9195484763085515
$row11[0,2] = 1.014008649095871
$row11[0,3] = 0.9794154155778751
$row11[0,4] = 1.037776720066717
$row11[0,5] = 0.9195484763085515
$row11[0,6] = 1.020246588260731
$row11[0,7] = 0.9799654006668854
$row11[0,8] = 1.010318185172167
$row11[0,9] = 0.9490305054455036
$row11[0,10] = 1.037776720066717
$row11[0,11] = 0.9667785627022113
$row11[0,12] = 0.9876873152622537
$row11[0,13] = 0.9887887425742876
$ws.Range("C11:P11").Value = $row11

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$row12 = New-Object 'object[,]' 1,14
$row12[0,0] = 1.03674519565764
$row12[0,1] = 0.9207194887095403
$row12[0,2] = 1.013876101259005
$row12[0,3] = 0.9796629449759158
$row12[0,4] = 1.03674519565764
$row12[0,5] = 0.9207194887095403
$row12[0,6] = 1.019840217749828
$row12[0,7] = 0.9802751382370257
$row12[0,8] = 1.009933091784984
$row12[0,9] = 0.9497855456935685
$row12[0,10] = 1.03674519565764
$row12[0,11] = 0.9672977949842727
$row12[0,12] = 0.9877509326505252
$row12[0,13] = 0.9888547155084384
$ws.Range("C12:P12").Value = $row12

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$row13 = New-Object 'object[,]' 1,14
$row13[0,0] = 1.037575532303584
$row13[0,1] = 0.9198451294928961
$row13[0,2] = 1.013944975176742
$row13[0,3] = 0.979475316320988
$row13[0,4] = 1.037575532303584
$row13[0,5] = 0.9198451294928961
$row13[0,6] = 1.020127931300284
$row13[0,7] = 0.9799929887866953
$row13[0,8] = 1.010224996543154
$row13[0,9] = 0.9491988078337664
$row13[0,10] = 1.037575532303584
$row13[0,11] = 0.9668950523348188
$row13[0,12] = 0.9877102383235524
$row13[0,13] = 0.9887982097197638
$ws.Range("C13:P13").Value = $row13

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$row14 = New-Object 'object[,]' 1,14
$row14[0,0] = 0.9917279999999995
$row14[0,1] = 0.8825720000000001
$row14[0,2] = 1.063179999999999
$row14[0,3] = 0.9718439999999998
$row14[0,4] = 0.9917279999999995
$row14[0,5] = 0.8825720000000001
$row14[0,6] = 1.035472000000002
$row14[0,7] = 1.006656
$row14[0,8] = 0.9960679999999991
$row14[0,9] = 0.9424160000000001
$row14[0,10] = 0.9917279999999995
$row14[0,11] = 0.9728759999999997
$row14[0,12] = 0.9773309999999997
$row14[0,13] = 0.986242
$ws.Range("C14:P14").Value = $row14

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$row15 = New-Object 'object[,]' 1,14
$row15[0,0] = 1.05
$row15[0,1] = 0.7441749999999996
$row15[0,2] = 1.12
$row15[0,3] = 0.94
$row15[0,4] = 1.05
$row15[0,5] = 0.7441749999999996
$row15[0,6] = 1.09
$row15[0,7] = 0.9997249999999999
$row15[0,8] = 1.02
$row15[0,9] = 0.86
$row15[0,10] = 1.05
$row15[0,11] = 0.9320874999999998
$row15[0,12] = 0.9635437499999999
$row15[0,13] = 0.9779875
$ws.Range("C15:P15").Value = $row15

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$row16 = New-Object 'object[,]' 1,14
$row16[0,0] = 1.027482800742397
$row16[0,1] = 0.8492699950080026
$row16[0,2] = 1.066959672217595
$row16[0,3] = 0.9633735897088004
$row16[0,4] = 1.027482800742397
$row16[0,5] = 0.8492699950080026
$row16[0,6] = 1.049439556607997
$row16[0,7] = 0.9962732007423972
$row16[0,8] = 1.009838548172802
$row16[0,9] = 0.9162895775744028
$row16[0,10] = 1.027482800742397
$row16[0,11] = 0.9581148336127987
$row16[0,12] = 0.9767715144191987
$row16[0,13] = 0.9848658675967994
$ws.Range("C16:P16").Value = $row16

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$row17 = New-Object 'object[,]' 1,14
$row17[0,0] = 0.9952493334421517
$row17[0,1] = 0.9950928326870866
$row17[0,2] = 0.9940714062568468
$row17[0,3] = 0.9943235584194052
$row17[0,4] = 0.9952493334421517
$row17[0,5] = 0.9950928326870866
$row17[0,6] = 0.9941345088924766
$row17[0,7] = 0.9946292140305567
$row17[0,8] = 0.9947240407802228
$row17[0,9] = 0.9940080201044336
$row17[0,10] = 0.9952204076170561
$row17[0,11] = 0.9945821194719666
$row17[0,12] = 0.9946842827013727
$row17[0,13] = 0.9945291143266475
$ws.Range("C17:P17").Value = $row17

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$row18 = New-Object 'object[,]' 1,14
$row18[0,0] = 0.9863588146071024
$row18[0,1] = 0.998589122701302
$row18[0,2] = 0.9975763477501917
$row18[0,3] = 0.9961446326860461
$row18[0,4] = 0.9863588146071024
$row18[0,5] = 0.998589122701302
$row18[0,6] = 0.993049372844873
$row18[0,7] = 0.9992302443422533
$row18[0,8] = 0.991994226048764
$row18[0,9] = 0.9968456381361711
$row18[0,10] = 0.9863588146071024
$row18[0,11] = 0.9980827352257469
$row18[0,12] = 0.9946672294361605
$row18[0,13] = 0.994973549889588
$ws.Range("C18:P18").Value = $row18

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$row19 = New-Object 'object[,]' 1,14
$row19[0,0] = 0.9806480779794212
$row19[0,1] = 1.016387545029663
$row19[0,2] = 0.9890765221817207
$row19[0,3] = 0.9995061223993642
$row19[0,4] = 0.9806480779794212
$row19[0,5] = 1.016387545029663
$row19[0,6] = 0.9857635524092959
$row19[0,7] = 0.998407628726652
$row19[0,8] = 0.9896564715863934
$row19[0,9] = 1.009251557959401
$row19[0,10] = 0.9806480779794212
$row19[0,11] = 1.002732033605692
$row19[0,12] = 0.9964045668975424
$row19[0,13] = 0.9960871847839889
$ws.Range("C19:P19").Value = $row19
